$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the employee identification data (RUN, Nombre, Apellido Paterno)
# for all data rows (2-13) - mantenedor usuario activo
$ws.Range("A2:A13").Value = "17459567-4"
$ws.Range("B2:B13").Value = "SERGIO"
$ws.Range("C2:C13").Value = "SOTO"

# Update the active cell selection to reflect where the user left off
$ws.Range("C16").Select()
